$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Capture the existing values before we shuffle rows around.
# (Use .Formula to read plain text reliably.)
$spotifyValue = $ws.Range("A9").Formula
$netflixValue = $ws.Range("A10").Formula

# Move the old "spotify" value into B9 and the old "netflix" value into C9,
# then put the new "streaming" category label in A9.
$ws.Range("B9").Value = $spotifyValue
$ws.Range("C9").Value = $netflixValue
$ws.Range("A9").Value = "streaming"

# The old row 10 (which held "netflix") is no longer needed on its own;
# remove it so everything below shifts up one row.
$ws.Rows.Item(10).Delete()

# Update the active selection to match the new layout.
$ws.Range("A10").Select()
